# إضافة حدث جديد في Card24 by admin at 2025-12-24 10:21:06
#
# Row 25's trailing placeholder cells (B..L and P..Q) were blank and get
# filled in with the sheet's usual "missing value" marker text "nan".
# A brand-new row 26 is appended below it holding the new service event
# (date / event / correction / serviced-by), with the same blank
# placeholder cells for the columns that don't apply.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Seed row 26 as a copy of row 25 first, while row 25 still has its
# original (blank) placeholder cells and a text-typed "24" in column A.
# That gives row 26 correctly-typed blank cells + a text "24" in A26
# without Excel re-interpreting it as a number.
$ws.Range("A25:Q25").Copy($ws.Range("A26:Q26"))

# --- Row 25: fill the previously-empty cells with "nan" ---------------
foreach ($col in @("B","C","D","E","F","G","H","I","J","K","L","P","Q")) {
    $ws.Range($col + "25").Value = "nan"
}

# --- Row 26: overwrite with the new event's actual data ----------------
$ws.Range("L26").Value = "22/12/2025"
$ws.Range("M26").Value = "كسر في شداد باب الدوفر الازرق"
$ws.Range("N26").Value = "تم تغيير الصموله الخاصة ب السداد"
$ws.Range("O26").Value = "ابراهيم/ناجي"
